$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column G width from 38 to 13 characters (ColumnWidth=12.17 -> stored width 13)
$ws.Columns.Item(7).ColumnWidth = 12.17

# Clear the "Recorded By" values (column G) for all filled session rows
$ws.Range("G2:G13").ClearContents()
$ws.Range("G15").ClearContents()
$ws.Range("G17:G30").ClearContents()
$ws.Range("G35").ClearContents()
$ws.Range("G37:G40").ClearContents()
$ws.Range("G42:G65").ClearContents()
$ws.Range("G67:G70").ClearContents()
$ws.Range("G77").ClearContents()
$ws.Range("G79:G80").ClearContents()
$ws.Range("G82").ClearContents()
$ws.Range("G84:G87").ClearContents()
$ws.Range("G92:G114").ClearContents()
$ws.Range("G116:G121").ClearContents()
$ws.Range("G127").ClearContents()
$ws.Range("G132:G154").ClearContents()
$ws.Range("G156:G165").ClearContents()
$ws.Range("G167:G172").ClearContents()
$ws.Range("G174:G175").ClearContents()
$ws.Range("G182:G205").ClearContents()
$ws.Range("G207:G212").ClearContents()
$ws.Range("G222:G230").ClearContents()
$ws.Range("G232:G235").ClearContents()
$ws.Range("G237:G256").ClearContents()
$ws.Range("G259:G261").ClearContents()
$ws.Range("G263").ClearContents()
$ws.Range("G270:G278").ClearContents()
$ws.Range("G282:G297").ClearContents()
$ws.Range("G299:G308").ClearContents()
$ws.Range("G310:G318").ClearContents()

# Data corrections that accompanied this processing pass
$ws.Range("H127").Value = "20/35"
$ws.Range("S18").Value = "68.8%"

